$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AA2").Value = 980
$ws.Range("AB2").Value = 14.5
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 12.5
$ws.Range("AH2").Value = 970
$ws.Range("AI2").Value = 980
$ws.Range("F2").Value = 3.95
$ws.Range("H2").Value = 2.06
$ws.Range("I2").Value = 2.18
$ws.Range("N2").Value = 3.4
$ws.Range("P2").Value = 1.81
$ws.Range("V2").Value = 1.84
$ws.Range("W2").Value = 1.29
$ws.Range("Y2").Value = 8.8
$ws.Range("Z2").Value = 15.5
$ws.Range("AD3").Value = 980
$ws.Range("AH3").Value = 980
$ws.Range("AK3").Value = 980
$ws.Range("AN3").Value = 10.5
$ws.Range("F3").Value = 1.46
$ws.Range("T3").Value = 1.91
$ws.Range("Y3").Value = 980
$ws.Range("AO4").Value = 55
$ws.Range("T4").Value = 1.83
$ws.Range("AA5").Value = 590
$ws.Range("AE5").Value = 280
$ws.Range("AJ5").Value = 9.6
$ws.Range("AO5").Value = 410
$ws.Range("G5").Value = 1.35
$ws.Range("H5").Value = 12.5
$ws.Range("P5").Value = 2
$ws.Range("T5").Value = 2.48
$ws.Range("U5").Value = 1.64
$ws.Range("W5").Value = 3.85
$ws.Range("K6").Value = 3.9
$ws.Range("T6").Value = 1.01
$ws.Range("K7").Value = 4.9
$ws.Range("Q7").Value = 1.73
$ws.Range("AF8").Value = 140
$ws.Range("AK8").Value = 130
$ws.Range("AM8").Value = 85
$ws.Range("AN8").Value = 80
$ws.Range("F8").Value = 11.5
$ws.Range("H8").Value = 1.27
$ws.Range("I8").Value = 1.28
$ws.Range("J8").Value = 7.6
$ws.Range("L8").Value = 1.15
$ws.Range("N8").Value = 11
$ws.Range("Q8").Value = 1.28
$ws.Range("V8").Value = 4.5
$ws.Range("W8").Value = 1.09
$ws.Range("AC9").Value = 9.8
$ws.Range("AF9").Value = 14
$ws.Range("N9").Value = 5.4
$ws.Range("Q9").Value = 1.63
$ws.Range("R9").Value = 1.6
$ws.Range("S9").Value = 2.56
$ws.Range("L10").Value = 1.28
$ws.Range("P10").Value = 2.48
$ws.Range("S10").Value = 2.58
$ws.Range("AD12").Value = 970
$ws.Range("AL12").Value = 970
$ws.Range("G12").Value = 1.54
$ws.Range("H12").Value = 6.8
$ws.Range("I12").Value = 7.8
$ws.Range("K12").Value = 5.3
$ws.Range("N12").Value = 1.03
$ws.Range("V12").Value = 1.15
$ws.Range("W12").Value = 2.88
$ws.Range("Y12").Value = 970
$ws.Range("AB13").Value = 9.6
$ws.Range("H13").Value = 26
$ws.Range("I13").Value = 27
$ws.Range("P13").Value = 2.72
$ws.Range("Q13").Value = 1.56
$ws.Range("Z13").Value = 370
$ws.Range("AE14").Value = 170
$ws.Range("I14").Value = 10
$ws.Range("AJ15").Value = 210
$ws.Range("AL15").Value = 90
$ws.Range("AO15").Value = 6.4
$ws.Range("F15").Value = 7
$ws.Range("H15").Value = 1.48
$ws.Range("I15").Value = 1.5
$ws.Range("J15").Value = 5.1
$ws.Range("P15").Value = 2.4
$ws.Range("Q15").Value = 1.59
$ws.Range("V15").Value = 3
$ws.Range("H16").Value = 1.04
$ws.Range("K16").Value = 980
$ws.Range("W16").Value = 3
